$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("route01")

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        if ($cell.Value2 -eq "*") {
            $cell.Value = "p"
        }
    }
}
